$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.169816493988037
$ws.Range("B1").Value = 4.367257118225098
$ws.Range("C1").Value = 8.122538566589355
$ws.Range("D1").Value = 8.375124931335449
$ws.Range("E1").Value = 5.524156093597412
